# Update workbook to add new expected-thanks columns and related test data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header cells (row 1)
$ws.Range("G1").Value = "expectedThanksHeader"
$ws.Range("H1").Value = "expectedThanksBody"

# New data cells (row 2)
$ws.Range("G2").Value = "Thank you for your order!"
$ws.Range("H2").Value = "Your order has been dispatched, and will arrive just as fast as the pony can get there!"

# Copy formatting from the existing header cells / text data cells so the
# new columns match the look of the rest of the table (reuses existing
# cell styles instead of creating new ones).
$ws.Range("A1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)

$ws.Range("D2").Copy()
$ws.Range("G2:H2").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Update the active selection to H1 to match the saved view state
$ws.Range("H1").Select()
